# Add a new "TOC" (table of contents) sheet as the first sheet in the
# workbook, listing every other sheet's logical table name and its
# (possibly truncated) tab name, then select it.

$wb = $excel.ActiveWorkbook

# Fix up a sheet name that picked up a better/longer name in this revision.
$wb.Worksheets.Item("TypeOfOfficerActivityCircumstan").Name = "OfficerActivityCircumstanceType"

# Insert the new sheet before the current first sheet so it becomes sheet #1.
$toc = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$toc.Name = "TOC"

# Column A = full/logical table name, Column B = actual worksheet tab name
# (tab names are limited to 31 chars, so some of these were truncated).
$tocData = @(
    ,("Table", "Tab")
    ,("SegmentActionTypeType", "SegmentActionTypeType")
    ,("ClearedExceptionallyType", "ClearedExceptionallyType")
    ,("UCROffenseCodeType", "UCROffenseCodeType")
    ,("OffenderSuspectedOfUsingType", "OffenderSuspectedOfUsingType")
    ,("LocationTypeType", "LocationTypeType")
    ,("MethodOfEntryType", "MethodOfEntryType")
    ,("TypeOfCriminalActivityType", "TypeOfCriminalActivityType")
    ,("TypeOfWeaponForceInvolvedType", "TypeOfWeaponForceInvolvedType")
    ,("BiasMotivationType", "BiasMotivationType")
    ,("TypePropertyLossEctType", "TypePropertyLossEctType")
    ,("PropertyDescriptionType", "PropertyDescriptionType")
    ,("SuspectedDrugTypeType", "SuspectedDrugTypeType")
    ,("TypeDrugMeasurementType", "TypeDrugMeasurementType")
    ,("TypeOfVictimType", "TypeOfVictimType")
    ,("OfficerActivityCircumstanceType", "OfficerActivityCircumstanceType")
    ,("OfficerAssignmentTypeType", "OfficerAssignmentTypeType")
    ,("SexOfPersonType", "SexOfPersonType")
    ,("RaceOfPersonType", "RaceOfPersonType")
    ,("EthnicityOfPersonType", "EthnicityOfPersonType")
    ,("ResidentStatusOfPersonType", "ResidentStatusOfPersonType")
    ,("AggravatedAssaultHomicideCircumstancesType", "AggravatedAssaultHomicideCircum")
    ,("AdditionalJustificationHomicideCircumstancesType", "AdditionalJustificationHomicideCi")
    ,("TypeInjuryType", "TypeInjuryType")
    ,("RelationshipsVictimToOffendersType", "RelationshipsVictimToOffendersT")
    ,("TypeOfArrestType", "TypeOfArrestType")
    ,("MultipleArresteeSegmentsIndicatorType", "MultipleArresteeSegmentsIndicat")
    ,("ArresteeWasArmedWithType", "ArresteeWasArmedWithType")
    ,("DispositionOfArresteeUnder18Type", "DispositionOfArresteeUnder18Typ")
)

$r = 1
foreach ($pair in $tocData) {
    $toc.Cells.Item($r, 1).Value = $pair[0]
    $toc.Cells.Item($r, 2).Value = $pair[1]
    $r = $r + 1
}

$toc.Columns.Item(1).ColumnWidth = 42
$toc.Columns.Item(2).ColumnWidth = 29

$toc.Range("B30").Select()
